$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# New rows 3 and 4 in column A need the same style as the existing A2 cell
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data in column B (shared strings "aaa" and "sf")
$ws.Range("B2").Value = "aaa"
$ws.Range("B3").Value = "sf"
